# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Several former sdmx/iaest *dimensions* (residencia-continente-nombre,
# inscripcion-provincia-nombre, inscripcion-municipio-estrato,
# residencia-area-nombre, sexo, inscripcion-comarca-nombre) are re-curated
# as *measures*. For each of those columns:
#   - Row 2 (field type) changes from the sdmx-dimension:refArea /
#     iaest-dimension:* identifier to the matching iaest-measure:* one.
#   - Row 3 (dim/medida marker) changes from "dim" to "medida".
#   - Row 4 (data type) changes from the skos:Concept / URI-* value to
#     "xsd:int".
#   - Row 5 (mapping file), where present, is removed since measures do
#     not carry a mapping workbook reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: residencia-continente-nombre
$ws.Range("B2").Value = "iaest-measure:residencia-continente-nombre"
$ws.Range("B3").Value = "medida"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("B5").Clear()

# Column G: inscripcion-provincia-nombre
$ws.Range("G2").Value = "iaest-measure:inscripcion-provincia-nombre"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"

# Column H: inscripcion-municipio-estrato
$ws.Range("H2").Value = "iaest-measure:inscripcion-municipio-estrato"
$ws.Range("H3").Value = "medida"
$ws.Range("H4").Value = "xsd:int"

# Column I: residencia-area-nombre
$ws.Range("I2").Value = "iaest-measure:residencia-area-nombre"
$ws.Range("I3").Value = "medida"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("I5").Clear()

# Column K: sexo
$ws.Range("K2").Value = "iaest-measure:sexo"
$ws.Range("K3").Value = "medida"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("K5").Clear()

# Column L: inscripcion-comarca-nombre
$ws.Range("L2").Value = "iaest-measure:inscripcion-comarca-nombre"
$ws.Range("L3").Value = "medida"
$ws.Range("L4").Value = "xsd:int"
